# "Add files via upload" — re-upload of "Lista de Archivos DATA.xlsx" with:
#   1) two "Fecha Actualización" (F-column) date corrections on the table rows
#      for Christian/Hector (row 35) and Edwin Hector (row 45): both moved to
#      the same updated date, serial 44154 (2020-11-19).
#   2) the "Situación" column (table column 4 / sheet column E) filtered down
#      to just "Trabajando", which hides every data row whose Situación isn't
#      "Trabajando" (rows 35 and 45 are the only survivors).
#   3) the resulting scroll position / active cell left where the user ended
#      up after filtering.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Fecha Actualización corrections ------------------------------------
$ws.Range("F35").Value = 44154
$ws.Range("F45").Value = 44154

# --- 2) Filter the table's "Situación" column down to "Trabajando" --------
$lo = $ws.ListObjects.Item(1)
$tableRange = $lo.Range
# Field 4 = column E ("Situación") within the table (B=1,C=2,D=3,E=4).
# Operator 7 = xlFilterValues -> emits <filters><filter val="..."/></filters>
# (a checkbox-style values filter) rather than a custom-compare filter.
$tableRange.AutoFilter(4, @("Trabajando"), 7)

# --- 3) Leave the view where it ended up after filtering ------------------
$ws.Range("F49").Select()
